$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(-0.424, -0.051, 0.845, 27, 60, -58),
    @(-0.075, -0.291, 0.656, -123, 47, 29),
    @(0.263, 0.08500000000000001, 0.35, -15, 25, 70),
    @(0.138, 0.179, 0.991, 90, -15, -121),
    @(0.231, 0.01, 0.795, 67, 15, -40),
    @(-0.227, -0.309, -0.508, 124, -2, 17),
    @(0.519, -0.447, 0.052, 46, 31, 33),
    @(-0.273, -0.369, 0.75, 175, 25, 170),
    @(0.001, -0.118, -0.12, 177, -29, -28),
    @(0.135, -0.238, -0.501, -84, -11, 157)
)

for ($r = 0; $r -lt $values.Length; $r++) {
    $row = $values[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}
